$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Teste de Fogo"
$ws.Range("A2").Value = 45602.95973272563
$ws.Range("B2").Value = "Teste de Fogo"
$ws.Range("C2").Value = 2424
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2500
$ws.Range("F2").Value = 5000

# Row 3 - "Garrafa de Água de Aço Inox" (unchanged text/code/unit price)
$ws.Range("A3").Value = 45602.96000383458
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 120

# Row 4 - "Livro de Receitas Vegetarianas"
$ws.Range("A4").Value = 45602.96056168922
$ws.Range("B4").Value = "Livro de Receitas Vegetarianas"
$ws.Range("C4").Value = 2345
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 59.9
$ws.Range("F4").Value = 59.9
